$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Adafruit IO pulled a new reading; append it as the next row directly
# below the existing data (the sheet currently ends at row 39).
$row = 40

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"

# "25" looks numeric, so force text formatting before assigning it,
# otherwise Excel would auto-convert it to a Number (the source data,
# like the rest of column C, is text). Reset the style afterwards so the
# cell doesn't end up with a lingering custom number format applied.
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "25"
$ws.Cells.Item($row, 3).Style = "Normal"

$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
